# Update the Omaha_B team-specific transition-probability matrix with
# recomputed probabilities after re-running the simulator with more games.
# Each cell below is a probability cell in the state-transition matrix
# (rows/cols A:S, header row 1); values were recalculated from updated
# simulated-game counts per state.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (state Af0)
$ws.Cells.Item(2, 2).Value = 0.1724137931034483
$ws.Cells.Item(2, 3).Value = 0.5862068965517241
$ws.Cells.Item(2, 10).Value = 0.03448275862068965
$ws.Cells.Item(2, 16).Value = 0.1413793103448276
$ws.Cells.Item(2, 19).Value = 0.06551724137931035

# Row 3 (state Af1)
$ws.Cells.Item(3, 3).Value = 0.04347826086956522
$ws.Cells.Item(3, 10).Value = 0.04347826086956522
$ws.Cells.Item(3, 16).Value = 0.7336956521739131
$ws.Cells.Item(3, 19).Value = 0.1793478260869565

# Row 4 (state Af2)
$ws.Cells.Item(4, 10).Value = 0.02777777777777778
$ws.Cells.Item(4, 16).Value = 0.6666666666666666
$ws.Cells.Item(4, 19).Value = 0.3055555555555556

# Row 5 (state Af3) - now has its first observed transition
$ws.Cells.Item(5, 16).Value = 1

# Row 6 (state Ai0)
$ws.Cells.Item(6, 2).Value = 0.04938271604938271
$ws.Cells.Item(6, 4).Value = 0.00411522633744856
$ws.Cells.Item(6, 6).Value = 0.06584362139917696
$ws.Cells.Item(6, 10).Value = 0.2510288065843622
$ws.Cells.Item(6, 15).Value = 0.01646090534979424
$ws.Cells.Item(6, 17).Value = 0.1769547325102881
$ws.Cells.Item(6, 18).Value = 0.07818930041152264
$ws.Cells.Item(6, 19).Value = 0.3580246913580247

# Row 7 (state Ai1)
$ws.Cells.Item(7, 2).Value = 0.09615384615384616
$ws.Cells.Item(7, 4).Value = 0.02307692307692308
$ws.Cells.Item(7, 5).Value = 0.003846153846153846
$ws.Cells.Item(7, 6).Value = 0.03461538461538462
$ws.Cells.Item(7, 10).Value = 0.1192307692307692
$ws.Cells.Item(7, 15).Value = 0.003846153846153846
$ws.Cells.Item(7, 17).Value = 0.1769230769230769
$ws.Cells.Item(7, 18).Value = 0.1
$ws.Cells.Item(7, 19).Value = 0.4423076923076923

# Row 8 (state Ai2)
$ws.Cells.Item(8, 2).Value = 0.1146881287726358
$ws.Cells.Item(8, 4).Value = 0.01006036217303823
$ws.Cells.Item(8, 6).Value = 0.07847082494969819
$ws.Cells.Item(8, 10).Value = 0.08853118712273642
$ws.Cells.Item(8, 15).Value = 0.01609657947686117
$ws.Cells.Item(8, 17).Value = 0.1670020120724346
$ws.Cells.Item(8, 18).Value = 0.1207243460764588
$ws.Cells.Item(8, 19).Value = 0.4044265593561368

# Row 9 (state Ai3)
$ws.Cells.Item(9, 2).Value = 0.1162790697674419
$ws.Cells.Item(9, 4).Value = 0.02325581395348837
$ws.Cells.Item(9, 6).Value = 0.05813953488372093
$ws.Cells.Item(9, 10).Value = 0.06976744186046512
$ws.Cells.Item(9, 15).Value = 0.02906976744186046
$ws.Cells.Item(9, 17).Value = 0.1918604651162791
$ws.Cells.Item(9, 18).Value = 0.1046511627906977
$ws.Cells.Item(9, 19).Value = 0.4069767441860465

# Row 10 (state Ar0)
$ws.Cells.Item(10, 2).Value = 0.09734513274336283
$ws.Cells.Item(10, 4).Value = 0.01609010458567981
$ws.Cells.Item(10, 5).Value = 0.0008045052292839903
$ws.Cells.Item(10, 6).Value = 0.07401448109412712
$ws.Cells.Item(10, 10).Value = 0.1263073209975865
$ws.Cells.Item(10, 15).Value = 0.02333065164923572
$ws.Cells.Item(10, 17).Value = 0.1922767497988737
$ws.Cells.Item(10, 18).Value = 0.09814963797264682
$ws.Cells.Item(10, 19).Value = 0.3716814159292036

# Row 11 (state Bf0)
$ws.Cells.Item(11, 7).Value = 0.1527777777777778
$ws.Cells.Item(11, 10).Value = 0.08333333333333333
$ws.Cells.Item(11, 11).Value = 0.1666666666666667
$ws.Cells.Item(11, 12).Value = 0.5944444444444444
$ws.Cells.Item(11, 19).Value = 0.002777777777777778

# Row 12 (state Bf1)
$ws.Cells.Item(12, 7).Value = 0.7465437788018433
$ws.Cells.Item(12, 10).Value = 0.1751152073732719
$ws.Cells.Item(12, 12).Value = 0.02764976958525346
$ws.Cells.Item(12, 19).Value = 0.05069124423963134

# Row 13 (state Bf2)
$ws.Cells.Item(13, 7).Value = 0.7758620689655172
$ws.Cells.Item(13, 10).Value = 0.1896551724137931
$ws.Cells.Item(13, 19).Value = 0.03448275862068965

# Row 15 (state Bi0)
$ws.Cells.Item(15, 6).Value = 0.0330188679245283
$ws.Cells.Item(15, 8).Value = 0.2122641509433962
$ws.Cells.Item(15, 9).Value = 0.02358490566037736
$ws.Cells.Item(15, 10).Value = 0.3584905660377358
$ws.Cells.Item(15, 11).Value = 0.09433962264150944
$ws.Cells.Item(15, 13).Value = 0.009433962264150943
$ws.Cells.Item(15, 15).Value = 0.04716981132075472
$ws.Cells.Item(15, 19).Value = 0.2216981132075472

# Row 16 (state Bi1)
$ws.Cells.Item(16, 6).Value = 0.01538461538461539
$ws.Cells.Item(16, 8).Value = 0.1230769230769231
$ws.Cells.Item(16, 9).Value = 0.07179487179487179
$ws.Cells.Item(16, 10).Value = 0.4461538461538462
$ws.Cells.Item(16, 11).Value = 0.1641025641025641
$ws.Cells.Item(16, 13).Value = 0.01538461538461539
$ws.Cells.Item(16, 14).Value = 0.005128205128205128
$ws.Cells.Item(16, 15).Value = 0.03589743589743589
$ws.Cells.Item(16, 19).Value = 0.1230769230769231

# Row 17 (state Bi2)
$ws.Cells.Item(17, 6).Value = 0.0136986301369863
$ws.Cells.Item(17, 8).Value = 0.1963470319634703
$ws.Cells.Item(17, 9).Value = 0.08447488584474885
$ws.Cells.Item(17, 10).Value = 0.4223744292237443
$ws.Cells.Item(17, 11).Value = 0.1118721461187215
$ws.Cells.Item(17, 13).Value = 0.0273972602739726
$ws.Cells.Item(17, 15).Value = 0.0639269406392694
$ws.Cells.Item(17, 19).Value = 0.07990867579908675

# Row 18 (state Bi3)
$ws.Cells.Item(18, 6).Value = 0.01659751037344398
$ws.Cells.Item(18, 8).Value = 0.1991701244813278
$ws.Cells.Item(18, 9).Value = 0.07468879668049792
$ws.Cells.Item(18, 10).Value = 0.4439834024896265
$ws.Cells.Item(18, 11).Value = 0.1078838174273859
$ws.Cells.Item(18, 13).Value = 0.02074688796680498
$ws.Cells.Item(18, 14).Value = 0.004149377593360996
$ws.Cells.Item(18, 15).Value = 0.06639004149377593
$ws.Cells.Item(18, 19).Value = 0.06639004149377593

# Row 19 (state Br0)
$ws.Cells.Item(19, 6).Value = 0.01923076923076923
$ws.Cells.Item(19, 8).Value = 0.2355769230769231
$ws.Cells.Item(19, 9).Value = 0.08253205128205128
$ws.Cells.Item(19, 10).Value = 0.3349358974358974
$ws.Cells.Item(19, 11).Value = 0.1290064102564103
$ws.Cells.Item(19, 13).Value = 0.03044871794871795
$ws.Cells.Item(19, 14).Value = 0.001602564102564103
$ws.Cells.Item(19, 15).Value = 0.06490384615384616
$ws.Cells.Item(19, 19).Value = 0.1017628205128205
